$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44299
$ws.Range("M2").Value2 = 80
$ws.Range("R2").Value2 = "Provincia de Santiago"

# Row 3
$ws.Range("D3").Value2 = 44299
$ws.Range("M3").Value2 = 75
$ws.Range("R3").Value2 = "Provincia de Santiago"

# Row 4
$ws.Range("D4").Value2 = 44301
$ws.Range("M4").Value2 = 100
$ws.Range("N4").Value2 = 14000
$ws.Range("O4").Value2 = 14000
$ws.Range("P4").Value2 = 14000
$ws.Range("S4").Value2 = 2000

# Row 5
$ws.Range("D5").Value2 = 44301
$ws.Range("M5").Value2 = 80
$ws.Range("N5").Value2 = 12000
$ws.Range("O5").Value2 = 12000
$ws.Range("P5").Value2 = 12000
$ws.Range("S5").Value2 = 1714

# Row 8
$ws.Range("D8").Value2 = 44300
$ws.Range("N8").Value2 = 15000
$ws.Range("O8").Value2 = 15000
$ws.Range("P8").Value2 = 15000
$ws.Range("S8").Value2 = 2143

# Row 9
$ws.Range("D9").Value2 = 44300

# Row 10
$ws.Range("D10").Value2 = 44302
$ws.Range("M10").Value2 = 50
$ws.Range("N10").Value2 = 15000
$ws.Range("O10").Value2 = 15000
$ws.Range("P10").Value2 = 15000
$ws.Range("S10").Value2 = 2143

# Row 11
$ws.Range("D11").Value2 = 44302
$ws.Range("M11").Value2 = 30
$ws.Range("N11").Value2 = 12000
$ws.Range("O11").Value2 = 12000
$ws.Range("P11").Value2 = 12000
$ws.Range("S11").Value2 = 1714

# Row 12
$ws.Range("D12").Value2 = 44322
$ws.Range("M12").Value2 = 45
$ws.Range("N12").Value2 = 12000
$ws.Range("O12").Value2 = 12000
$ws.Range("P12").Value2 = 12000
$ws.Range("R12").Value2 = "Región Metropolitana"
$ws.Range("S12").Value2 = 1714

# Row 13
$ws.Range("D13").Value2 = 44322
$ws.Range("M13").Value2 = 80
$ws.Range("N13").Value2 = 8000
$ws.Range("O13").Value2 = 8000
$ws.Range("P13").Value2 = 8000
$ws.Range("R13").Value2 = "Región Metropolitana"
$ws.Range("S13").Value2 = 1143

# Row 14
$ws.Range("D14").Value2 = 44320
$ws.Range("M14").Value2 = 20
$ws.Range("N14").Value2 = 12000
$ws.Range("O14").Value2 = 12000
$ws.Range("P14").Value2 = 12000
$ws.Range("S14").Value2 = 1714

# Row 15
$ws.Range("D15").Value2 = 44320
$ws.Range("N15").Value2 = 8000
$ws.Range("O15").Value2 = 8000
$ws.Range("P15").Value2 = 8000
$ws.Range("S15").Value2 = 1143
